# ---------------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table changes its table-style (tableStyleId GUID) from
#    {21C4C9CB-3600-4E12-A68D-4F5BC4D24227} to
#    {46B39096-7CFC-4C9B-830D-8AC70DFAEBFC}.
#
# 2) The design theme ("Office Theme" blue palette) is (re)applied to the
#    deck's slide master / theme part, replacing the "Integral" (Red Violet)
#    palette that was in use. This is the swap that moves the "Office Theme"
#    colour values onto the theme part that slideMaster1.xml / the
#    presentation actually renders with.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 --------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{46B39096-7CFC-4C9B-830D-8AC70DFAEBFC}")

# --- 2) Swap the theme colour scheme back to the default "Office Theme" ---
$theme = $p.Designs.Item(1).SlideMaster.Theme

# Cosmetic (best-effort) naming - harmless if the host ignores it.
$theme.Name = "Office Theme"

$colorScheme = $theme.ThemeColorScheme

# Index : role      : Office Theme RGB
#   1   : Dark 1     : 000000
#   2   : Light 1    : FFFFFF
#   3   : Dark 2     : 44546A
#   4   : Light 2    : E7E6E6
#   5   : Accent 1   : 5B9BD5
#   6   : Accent 2   : ED7D31
#   7   : Accent 3   : A5A5A5
#   8   : Accent 4   : FFC000
#   9   : Accent 5   : 4472C4
#  10   : Accent 6   : 70AD47
#  11   : Hyperlink  : 0563C1
#  12   : Followed   : 954F72
$colorScheme.Item(1).RGB  = 0x000000
$colorScheme.Item(2).RGB  = 0xFFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444
$colorScheme.Item(4).RGB  = 0xE6E6E7
$colorScheme.Item(5).RGB  = 0xD59B5B
$colorScheme.Item(6).RGB  = 0x317DED
$colorScheme.Item(7).RGB  = 0xA5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF
$colorScheme.Item(9).RGB  = 0xC47244
$colorScheme.Item(10).RGB = 0x47AD70
$colorScheme.Item(11).RGB = 0xC16305
$colorScheme.Item(12).RGB = 0x724F95
